$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "21.896.24"
$ws.Range("E2").Value = "  +6.61%  "
$ws.Range("D3").Value = "1.576.33"
$ws.Range("E3").Value = "  +6.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9880"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.81%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3696"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3289"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.56%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "41.98"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.61%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.145"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.96%  "
$ws.Range("E11").Value = "  +6.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9985"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +11.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.874"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.548"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9881"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001076"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.24%  "
$ws.Range("D18").Value = "1.575.82"
$ws.Range("E18").Value = "  +6.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06402"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "75.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.08%  "
$ws.Range("E21").Value = "  +11.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.879"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.80%  "
$ws.Range("E23").Value = "  +5.67%  "
$ws.Range("D24").Value = "21.914.32"
$ws.Range("E24").Value = "  +6.52%  "
$ws.Range("E25").Value = "  +5.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.426"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +13.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.50%  "
$ws.Range("D29").Value = "1.748.34"
$ws.Range("E29").Value = "  +6.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.163"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9255"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +13.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.645"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +13.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08227"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.632"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.82%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.703"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +13.38%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +13.28%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.204"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06170"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.241"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02186"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2007"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9874"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5839"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.650"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5675"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.943"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06799"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.90%  "
